$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the previous B (dbExcel)
# and C (WebExcel) columns one slot to the right (to C and D), matching the
# diff's column remap: old B->C, old C->D.
$ws.Columns("B:B").Insert()

# New header for the inserted "StatQuery" column.
$ws.Range("B1").Value = "StatQuery"

# New stat-bar Neo4j query text placed under the new header, row 2.
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Vizsla']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Give the new B2 cell the same wrap-text style used by A2 (style index 1).
$ws.Range("B2").WrapText = $true

# Re-apply the original column widths to the (now shifted) columns, and
# give the new column B the same width as column A (75.81640625 / 70.26953125
# / 28.54296875 characters, as in the source workbook).
$ws.Columns("A:A").ColumnWidth = 74.98307291666667
$ws.Columns("B:B").ColumnWidth = 74.98307291666667
$ws.Columns("C:C").ColumnWidth = 69.43619791666667
$ws.Columns("D:D").ColumnWidth = 27.709635416666668

# The saved view no longer scrolls to row 2 - restore the top-left cell/selection.
$ws.Application.GoTo($ws.Range("A1")) | Out-Null
$ws.Range("A2").Select() | Out-Null
